$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Bvals = @(0.1423106037013184,0.1330144244095095,0.1273773673674725,0.1250981766461052,0.1247208074859572,0.1273465565114407,0.1390906624479697,0.1626778301107947,0.1803420342742896,0.1884497011099313,0.1915301156894884,0.190866241360979,0.1887029247956775,0.1873791568901453,0.1798136173794376,0.1751907681310314,0.1725386261944806,0.1716418289157389,0.1756821758123124,0.1893380672599108,0.1983224626450237,0.1935219306116522,0.1754599929582525,0.1562377461785047)
$Dvals = @(0.1697448963614363,0.1671367038718614,0.1655945231664901,0.1649810464549546,0.164880085687976,0.1655861888770076,0.1688333348001407,0.1756684330423042,0.1809721415735481,0.1834455681805451,0.1843908615254435,0.1841868916996248,0.1835231648906301,0.1831177386883098,0.1808117132666354,0.1794125460931042,0.1786135039940575,0.1783439470728609,0.1795608981769448,0.18371788309571,0.1864851579397566,0.1850036200878122,0.1794938115062905,0.1737696258197872)
$Evals = @(0.1564731755020787,0.1560147298383896,0.1558024348625224,0.1557333510515306,0.15572293320988,0.1558014325694614,0.1563007594761672,0.1578278911207462,0.1592828238758592,0.1600167961721688,0.1603050785956057,0.1602425322405985,0.1600403062062661,0.1599177830845235,0.1592363058983999,0.1588366930520202,0.1586136358445991,0.1585392794429623,0.158878530027625,0.1600994244038105,0.1609576240263983,0.1604940805836641,0.1588595947045235,0.1573562129349853)
$Fvals = @(1.559342287038866,1.56552637913456,1.570194383556817,1.572315649587104,1.572681112539335,1.570222104893162,1.561293800236456,1.550696923316536,1.547128097068651,1.546421002019798,1.546285054590541,1.546308470211983,1.546407175568206,1.546484802697947,1.54719275073515,1.547861789678819,1.548332855064842,1.548507162341721,1.547781642303974,1.546374605703875,1.546223357326255,1.546233766199464,1.547817607721214,1.552823396186376)
$Gvals = @(0.002474408482442186,0.002477202841472538,0.002479010731127067,0.002479770700129648,0.002479898297944929,0.002479020886256443,0.002475352895005455,0.00246888791470906,0.002464577427098863,0.002462710927446332,0.002462017630495969,0.002462166344861069,0.002462653619178386,0.002462953845795082,0.002464701300519467,0.00246579743077931,0.002466436781727754,0.002466654783009268,0.002465679826725906,0.002462510128789203,0.002460517242311928,0.002461573703572273,0.002465732966990956,0.002470559392974996)
$Jvals = @(0.1845980772993379,0.1853464863728931,0.1858982256128989,0.1861462592172174,0.1861888461790073,0.1859014767535072,0.1848369942782213,0.1834810108745302,0.1829306082935105,0.1827770231735073,0.182732780395618,0.1827416899667611,0.1827731043992458,0.1827941589188029,0.1829425930042419,0.1830584424312605,0.183134186988859,0.1831613976998838,0.1830451671227991,0.1827634995475691,0.1826605287429217,0.1827080649978043,0.1830511404162607,0.1834810108745302)
$Kvals = @(0.6267408329052273,0.5517257592008491,0.5055742526283495,0.4867451006184069,0.4836172352496533,0.505320403890579,0.6008953437470268,0.7875496487299358,0.9241786514892851,0.986217799202592,1.009693120684744,1.00463808932858,0.9881494864850708,0.9780474205551002,0.9201218530370454,0.884556353426774,0.8640893502730194,0.8571577776641561,0.8883434741079839,0.9929930744862077,1.061284770764416,1.024845998689102,0.8866313788232105,0.7371408604769556)
$Mvals = @(0.2682518545108366,0.2500296013090733,0.2389000367067027,0.2343797446631086,0.2336300732559806,0.2388390129895797,0.2619567355251604,0.3077491461840722,0.3416635093013412,0.3571491577134864,0.363021270552224,0.361756253162838,0.3576321005737668,0.3551069759143601,0.3406526253387199,0.3317999606318338,0.3267136024953885,0.3249923955032443,0.3327417788184945,0.3588432483683803,0.3759487650334492,0.3668150490415911,0.332315972962995,0.2953129141995561)
$Ovals = @(3.849845266364355,3.874377125917931,3.891704648677916,3.899335049226693,3.900636448868909,3.89180525018196,3.857833641294746,3.80919846068241,3.78444960191635,3.775580300422263,3.772565608057221,3.773199574157388,3.775325385075462,3.776672305171871,3.785077337923212,3.790845719589441,3.794388389315174,3.795626483089848,3.790208391753595,3.774691646034228,3.766555354381723,3.770714284348173,3.79049582261905,3.820428237442314)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $Bvals[$i]
    $ws.Cells.Item($row, 4).Value = $Dvals[$i]
    $ws.Cells.Item($row, 5).Value = $Evals[$i]
    $ws.Cells.Item($row, 6).Value = $Fvals[$i]
    $ws.Cells.Item($row, 7).Value = $Gvals[$i]
    $ws.Cells.Item($row, 10).Value = $Jvals[$i]
    $ws.Cells.Item($row, 11).Value = $Kvals[$i]
    $ws.Cells.Item($row, 13).Value = $Mvals[$i]
    $ws.Cells.Item($row, 15).Value = $Ovals[$i]
}
